# edit.ps1 - applies the tracked-change-style rewrite of the recommendation
# letter body text, paragraph structure and signature block.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# ---------------------------------------------------------------------
# Paragraph 3 (intro) - text-only changes, paragraph formatting unchanged
# ---------------------------------------------------------------------
Replace-Text "I would like to recommend Milan for admission to the graduate program at your university. I " "It is my pleasure to write this letter of recommendation for milan. I have known them for about 6 "
Replace-Text "have known for about 5 now as an undergraduate student in BCT Engineering. Moreover, I was " "years as an undergraduate student in electronics and computer. I have taught them in several "
Replace-Text "the supervisor for his final year project. I taught him and ." "courses. I have also mentored and evaluated their project work."

# ---------------------------------------------------------------------
# Paragraph 4 (academic performance) - text changes + new trailing sentence
# ---------------------------------------------------------------------
Replace-Text "I recall as a student. He maintained excellent academic performance throughout his " "They has consistently demonstrated sharp intellect and a strong work ethic. Academically, they "
Replace-Text "undergraduate ranking among the students of his batch. As an instructor and his supervisor too " "has maintained excellent performance throughout the undergraduate years, showing strength "
Replace-Text "I had enough opportunity is observe his capabilities closely." "in both theoretical understanding and practical problem-solving. They is also a versatile and quick learner."

# paragraph spacing-before changes from 326 -> 0 twips (16.3 -> 0 pt)
$pAcademic = $d.Paragraphs.Item(4)
$pAcademic.Format.SpaceBefore = 0

# ---------------------------------------------------------------------
# Insert two brand-new paragraphs right after the academic paragraph
# ---------------------------------------------------------------------
$rng = $pAcademic.Range
$rng.InsertParagraphAfter()
$pProject1 = $d.Paragraphs.Item(5)
$pProject1.Range.Text = 'I particularly recall the third-year project titled "Recommendation Generator", which required extensive study in areas like AI and algorithms. The project demonstrated their self-learning ability and research interest.'
$pProject1.Format.SpaceBefore = 0
$pProject1.Format.LineSpacing = 16.4
$pProject1.Format.LineSpacingRule = 4
$pProject1.Format.LeftIndent = 1.7
$pProject1.Format.RightIndent = 1.7
$pProject1.Format.Alignment = 3

$rng2 = $pProject1.Range
$rng2.InsertParagraphAfter()
$pProject2 = $d.Paragraphs.Item(6)
$pProject2.Range.Text = 'As part of a DBMS course, they also led a team that built a system titled "No", where they served as the team lead. I had the opportunity to observe their leadership and teamwork closely in both projects.'
$pProject2.Format.SpaceBefore = 0
$pProject2.Format.LineSpacing = 16.4
$pProject2.Format.LineSpacingRule = 4
$pProject2.Format.LeftIndent = 1.7
$pProject2.Format.RightIndent = 1.7
$pProject2.Format.Alignment = 3

# ---------------------------------------------------------------------
# Old paragraph "I was the supervisor in his project titled ." becomes the
# short "Beyond the classroom..." extracurricular paragraph.
# ---------------------------------------------------------------------
Replace-Text "I was the supervisor in his project titled . I was quite impressed by his hardworking nature and " "Beyond the classroom, they actively participated in extracurricular activities such as "
Replace-Text "learning capability. In fact, he was also able to publish a paper on ." "Hackathon."

$pExtra = $d.Paragraphs.Item(7)
$pExtra.Format.LineSpacing = 16.3
$pExtra.Format.LineSpacingRule = 4
$pExtra.Format.SpaceBefore = 0

# ---------------------------------------------------------------------
# Delete the three obsolete paragraphs that followed:
#   "I have noted his presentation skills..."
#   "I appreciate his technical and professional skills..."
#   "I am quite confident that 's skills..."
# ---------------------------------------------------------------------
$d.Paragraphs.Item(8).Range.Delete() | Out-Null
$d.Paragraphs.Item(8).Range.Delete() | Out-Null
$d.Paragraphs.Item(8).Range.Delete() | Out-Null

# ---------------------------------------------------------------------
# The old signature paragraph ("DSB, " / "CIT Admin, " / ...) becomes the
# "On a personal level..." paragraph; its first run is rewritten and a new
# second run/sentence is appended.
# ---------------------------------------------------------------------
Replace-Text "DSB, " "On a personal level, they is polite, easy-going, and friendly. It has been a pleasure to teach them, and I have always found them to be respectful and open to feedback."

$pPersonal = $d.Paragraphs.Item(8)
$pPersonal.Format.SpaceBefore = 0
$pPersonal.Format.RightIndent = 0
$pPersonal.Format.LeftIndent = 1.7
$pPersonal.Format.LineSpacing = 16.4
$pPersonal.Format.LineSpacingRule = 4
$pPersonal.Format.Alignment = 0

# ---------------------------------------------------------------------
# Insert the remaining new closing paragraphs before the signature block.
# ---------------------------------------------------------------------
$rngP = $pPersonal.Range
$rngP.InsertParagraphAfter()
$pHarvard = $d.Paragraphs.Item(9)
$pHarvard.Range.Text = "I am confident that if given the opportunity to pursue the computer program at Harvard, they will excel and contribute meaningfully."
$pHarvard.Format.SpaceBefore = 0
$pHarvard.Format.LineSpacing = 16.3
$pHarvard.Format.LineSpacingRule = 4
$pHarvard.Format.LeftIndent = 1.7
$pHarvard.Format.RightIndent = 0
$pHarvard.Format.Alignment = 0

$rngH = $pHarvard.Range
$rngH.InsertParagraphAfter()
$pContact = $d.Paragraphs.Item(10)
$pContact.Range.Text = "Please feel free to contact me at dsb@gmail.com if you need any further information."
$pContact.Format.SpaceBefore = 1.3
$pContact.Format.LineSpacing = 15.1
$pContact.Format.LineSpacingRule = 4
$pContact.Format.LeftIndent = 1.7
$pContact.Format.RightIndent = 0
$pContact.Format.Alignment = 0

# ---------------------------------------------------------------------
# Signature block: text + indentation changes. The old "DSB, " run was
# already consumed above, so the remaining runs are updated in place and
# the paragraph's right indent widens (2880 -> 5040 twips).
# ---------------------------------------------------------------------
Replace-Text "CIT Admin, " "DSB, Cit Admin "
Replace-Text "Department of Electronics and Computer Engineering " "Department of Electronics and Computer "
Replace-Text "Pulchowk Campus, Institute of Engineering, Tribhuvan University " "Institute of Engineering, Pulchowk Campus "
Replace-Text "Phone: 9876543210 " ""
Replace-Text "Mail: dsb@gmail.com" "dsb@gmail.com"

$pSig = $d.Paragraphs.Item(11)
$pSig.Format.SpaceBefore = 0
$pSig.Format.RightIndent = 252

Write-Host "Done"
